$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("L2").Value = 1.29
$ws.Range("M2").Value = 3.75
$ws.Range("N2").Value = 1.97
$ws.Range("O2").Value = 1.93

# Row 3
$ws.Range("T3").Value = 10
$ws.Range("Z3").Value = 17
$ws.Range("AA3").Value = 9.5
$ws.Range("AH3").Value = 81

# Row 4
$ws.Range("G4").Value = 2.88
$ws.Range("I4").Value = 2.4
$ws.Range("K4").Value = 9.5
$ws.Range("T4").Value = 9
$ws.Range("AD4").Value = 251
$ws.Range("AE4").Value = 8

# Row 8
$ws.Range("J8").Value = 1.06
$ws.Range("K8").Value = 10
$ws.Range("R8").Value = 2.1
$ws.Range("S8").Value = 1.67
$ws.Range("T8").Value = 15
$ws.Range("W8").Value = 67
$ws.Range("AA8").Value = 8.5
$ws.Range("AB8").Value = 21
$ws.Range("AE8").Value = 6
$ws.Range("AH8").Value = 9.5

# Row 9
$ws.Range("K9").Value = 9.5
$ws.Range("R9").Value = 1.8
$ws.Range("S9").Value = 1.91
$ws.Range("W9").Value = 51
$ws.Range("Z9").Value = 9.5
$ws.Range("AC9").Value = 51
$ws.Range("AF9").Value = 8.5

# Row 13
$ws.Range("G13").Value = 1.6
$ws.Range("I13").Value = 4.75
$ws.Range("W13").Value = 13

# Row 14
$ws.Range("G14").Value = 2.5
$ws.Range("I14").Value = 3
$ws.Range("X14").Value = 23
$ws.Range("Y14").Value = 41
$ws.Range("AE14").Value = 7.5
$ws.Range("AF14").Value = 13
$ws.Range("AG14").Value = 11
$ws.Range("AH14").Value = 29
$ws.Range("AI14").Value = 26

# Row 17
$ws.Range("G17").Value = 1.65
$ws.Range("H17").Value = 3.35
$ws.Range("I17").Value = 4.8
$ws.Range("N17").Value = 1.9
$ws.Range("O17").Value = 1.72
$ws.Range("U17").Value = 6.3
$ws.Range("W17").Value = 10.5
$ws.Range("X17").Value = 11.5
$ws.Range("Z17").Value = 9
$ws.Range("AA17").Value = 5.8
$ws.Range("AB17").Value = 12.5
$ws.Range("AC17").Value = 55
$ws.Range("AD17").Value = 400
$ws.Range("AE17").Value = 10.75
$ws.Range("AF17").Value = 23
$ws.Range("AG17").Value = 12.5
$ws.Range("AH17").Value = 65
$ws.Range("AI17").Value = 37
$ws.Range("AJ17").Value = 37

# Row 19
$ws.Range("H19").Value = 7.5
$ws.Range("I19").Value = 17.5
$ws.Range("T19").Value = 8.5
$ws.Range("U19").Value = 5.6
$ws.Range("W19").Value = 5.2
$ws.Range("X19").Value = 10
$ws.Range("Y19").Value = 35
$ws.Range("Z19").Value = 19.5
$ws.Range("AA19").Value = 16
$ws.Range("AB19").Value = 37
$ws.Range("AC19").Value = 175
$ws.Range("AF19").Value = 120
$ws.Range("AJ19").Value = 175

# Row 24
$ws.Range("H24").Value = 6.25
$ws.Range("I24").Value = 10
$ws.Range("U24").Value = 7
$ws.Range("V24").Value = 10
$ws.Range("AA24").Value = 13

# Row 25
$ws.Range("J25").Value = 1.04
$ws.Range("K25").Value = 9

# Row 28
$ws.Range("G28").Value = 1.55
$ws.Range("H28").Value = 4.1
$ws.Range("I28").Value = 5.5
$ws.Range("M28").Value = 4.05
$ws.Range("R28").Value = 1.72
$ws.Range("S28").Value = 2
$ws.Range("T28").Value = 7.4
$ws.Range("U28").Value = 8.25
$ws.Range("W28").Value = 12
$ws.Range("AA28").Value = 8.5
$ws.Range("AB28").Value = 16.5
$ws.Range("AC28").Value = 70
$ws.Range("AE28").Value = 16
$ws.Range("AF28").Value = 40
$ws.Range("AG28").Value = 18
$ws.Range("AH28").Value = 120
$ws.Range("AI28").Value = 55
$ws.Range("AJ28").Value = 50

# Row 30
$ws.Range("G30").Value = 1.25
$ws.Range("H30").Value = 5.4
$ws.Range("I30").Value = 10.5
$ws.Range("J30").Value = 1.03
$ws.Range("K30").Value = 9.5
$ws.Range("L30").Value = 1.14
$ws.Range("M30").Value = 4.85
$ws.Range("N30").Value = 1.45
$ws.Range("O30").Value = 2.55
$ws.Range("P30").Value = 1.25
$ws.Range("Q30").Value = 3.55
$ws.Range("R30").Value = 1.88
$ws.Range("S30").Value = 1.83
$ws.Range("T30").Value = 8.75
$ws.Range("U30").Value = 6.9
$ws.Range("V30").Value = 8.75
$ws.Range("W30").Value = 7.9
$ws.Range("X30").Value = 10.25
$ws.Range("Y30").Value = 24
$ws.Range("Z30").Value = 9.5
$ws.Range("AA30").Value = 11.25
$ws.Range("AB30").Value = 21
$ws.Range("AC30").Value = 80
$ws.Range("AD30").Value = 500
$ws.Range("AE30").Value = 32
$ws.Range("AF30").Value = 90
$ws.Range("AJ30").Value = 80

# Row 31
$ws.Range("G31").Value = 3.8
$ws.Range("H31").Value = 3.85
$ws.Range("I31").Value = 1.8
$ws.Range("K31").Value = 8.25
$ws.Range("L31").Value = 1.23
$ws.Range("M31").Value = 3.75
$ws.Range("P31").Value = 1.35
$ws.Range("Q31").Value = 2.95
$ws.Range("R31").Value = 1.7
$ws.Range("S31").Value = 2.05
$ws.Range("T31").Value = 12.5
$ws.Range("Z31").Value = 8.25
$ws.Range("AA31").Value = 7.5
$ws.Range("AE31").Value = 8.25
$ws.Range("AH31").Value = 15
$ws.Range("AI31").Value = 13.5
